$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1829735462431551
$ws.Range("C2").Value = 0.6694626356775797
$ws.Range("D2").Value = 0.8372369824267251
$ws.Range("E2").Value = 0.9150065477507389
$ws.Range("F2").Value = 0.905446197782048
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.1263938744343593
$ws.Range("C3").Value = 0.8614631615508216
$ws.Range("D3").Value = 1.473077372104747
$ws.Range("E3").Value = 1.213703988666408
$ws.Range("F3").Value = 1.219359969539984
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = -0.09485759795785621
$ws.Range("C4").Value = 0.6206824703513352
$ws.Range("D4").Value = 0.6573762981246858
$ws.Range("E4").Value = 0.8107874555792571
$ws.Range("F4").Value = 0.813563898042168
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.05113506990239874
$ws.Range("C5").Value = 0.7519750174941408
$ws.Range("D5").Value = 1.034754511348446
$ws.Range("E5").Value = 1.017228839223725
$ws.Range("F5").Value = 1.026693788557885
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = -0.03772088114894339
$ws.Range("C6").Value = 0.5996523225882197
$ws.Range("D6").Value = 0.6704444970345537
$ws.Range("E6").Value = 0.8188067519473405
$ws.Range("F6").Value = 0.8267802493317774
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.1077066685141252
$ws.Range("C7").Value = 0.790586417127402
$ws.Range("D7").Value = 1.154731043345991
$ws.Range("E7").Value = 1.074584125764936
$ws.Range("F7").Value = 1.080987147395736
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = -0.03868616839471213
$ws.Range("C8").Value = 0.5498893294951147
$ws.Range("D8").Value = 0.5944205642755658
$ws.Range("E8").Value = 0.7709867471465159
$ws.Range("F8").Value = 0.7787165651551944
$ws.Range("G8").Value = 45

$ws.Range("B9").Value = 0.1749977985370214
$ws.Range("C9").Value = 0.7501568359084184
$ws.Range("D9").Value = 1.058248496194531
$ws.Range("E9").Value = 1.028712056988996
$ws.Range("F9").Value = 1.025437714082251
$ws.Range("G9").Value = 44

$ws.Range("B10").Value = 0.02191909566903653
$ws.Range("C10").Value = 0.5542941723622423
$ws.Range("D10").Value = 0.633576040990445
$ws.Range("E10").Value = 0.7959748997238826
$ws.Range("F10").Value = 0.8050896216323691
$ws.Range("G10").Value = 43

$ws.Range("B11").Value = 0.2576676093386086
$ws.Range("C11").Value = 0.788856672248688
$ws.Range("D11").Value = 1.23948852930275
$ws.Range("E11").Value = 1.113323191756441
$ws.Range("F11").Value = 1.096224442490078
$ws.Range("G11").Value = 42
